# report-checklist_WERFEN_MODULAB.xlsx edits
# "correzione workflowInstanceId test 1, 2, 3, 4, 5" +
# "correzione errori test 4" +
# "correzione test 56 - eliminazione test in quanto non riproducibile"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# --- Test 1 (row 8): strip the "^^^^urn:ihe:iti:xdw:2013:workflowInstanceId" suffix from WORKFLOWINSTANCEID
$ws.Range("I8").Value = "2.16.840.1.113883.2.9.2.100.4480c3c0f735720f2cb9c63251201c888e6bf05275b53d7ba161280b2eaa9a92.96872754e5"

# --- Test 2 (row 9): strip the suffix
$ws.Range("I9").Value = "2.16.840.1.113883.2.9.2.100.49dc14dd85651002abad9bc6bea92845e90f8d1e6c87cffa35a7e1f8fb7766dd.f0ec2f6e4e"

# --- Test 3 (row 10): strip the suffix
$ws.Range("I10").Value = "2.16.840.1.113883.2.9.2.100.52325a07a1fdcae7cd4eb95334f4d47a7260cc8851e7741f13d29fb5e22e2f5e.fa6c0cd8c7"

# --- Test 4 (row 11): correct the execution data entirely (new date/timestamp/traceid/workflowid)
$ws.Range("F11").Value = "03/24/2023"
$ws.Range("G11").Value = "2023-03-24T11:19:18Z"
$ws.Range("H11").Value = "215f6c833e6e2201"
$ws.Range("I11").Value = "2.16.840.1.113883.2.9.2.30.4d916594b2eb4957cde85fd788c191b276fafe65c3135df4f6398dcfaa38fc91.474a921572"

# --- Test 5 (row 12): strip the suffix
$ws.Range("I12").Value = "2.16.840.1.113883.2.9.2.100.6b785ecfa5d59e3a14f59fd2d37b696a45d71bcbb41c0326fb2bd24e40ec4331.72b3a92686"

# --- Test 56 (row 20): eliminated as non-reproducible -> clear execution data, mark not applicable
$ws.Range("F20").Value = ""
$ws.Range("G20").Value = ""
$ws.Range("H20").Value = ""
$ws.Range("I20").Value = ""
$ws.Range("J20").Value = "NO"
$ws.Range("K20").Value = "Campo sempre presente e valorizzato in maniera corretta"
$ws.Range("L20").Value = ""
$ws.Range("M20").Value = ""
$ws.Range("O20").Value = ""
$ws.Range("P20").Value = ""

# Restore view state: top-left cell + selection like the saved file
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("P20").Select()
